# weekly-activities.xlsx update
# - Tuesday (B3) task renamed from "Aguia - 1 hora" to "Rever front - 1 hora"
# - Wednesday (row 4): "Elaborar Contrato - 1 hora" (C4) is dropped; the
#   remaining "Marcar reuniao - 1 hora" task shifts from D4 into C4, and the
#   now-trailing D4 cell is removed.
# - Friday (row 6): the two extra tasks "Reuniao Sinavez - 1 hora" (C6) and
#   "Implementar front - 1 hora" (D6) are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename Tuesday's first task.
$ws.Range("B3").Value = "Rever front - 1 hora"

# Wednesday: drop "Elaborar Contrato - 1 hora" and shift "Marcar reuniao - 1 hora" left.
$ws.Range("C4").Value = "Marcar reuniao - 1 hora"
$ws.Range("D4").Clear()

# Friday: remove the last two tasks of the week.
$ws.Range("C6").Clear()
$ws.Range("D6").Clear()
